$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-33 down to 13-34.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112039
$ws.Cells.Item(12, 7).Value = "Ciboulette"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = 2000
$ws.Cells.Item(12, 14).Value = "`$/docena de atados"
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 667
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = "Hortaliza"
